# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) previously held a raw strikeout-count style
# statistic; it is regenerated here so that each row reflects the
# recalculated value actually written by the save routine.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new value for column G ("K")
$kValues = @{
    2  = 6
    3  = 4
    4  = 3
    5  = 0
    6  = 3
    7  = 4
    8  = 2
    9  = 4
    10 = 4
    11 = 3
    12 = 2
    13 = 4
    14 = 4
    15 = 5
    16 = 2
    17 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
